$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.675.59"
$ws.Range("E2").Value = "  +6.22%  "
$ws.Range("D3").Value = "2.042.62"
$ws.Range("E3").Value = "  +3.38%  "
$ws.Range("D5").Value = "'251.05"
$ws.Range("E5").Value = "  +4.35%  "
$ws.Range("E6").Value = "  +2.44%  "
$ws.Range("D7").Value = "'65.84"
$ws.Range("E7").Value = "  +16.52%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'59.55"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'0.374"
$ws.Range("E10").Value = "  +4.66%  "
$ws.Range("D11").Value = "'0.0754"
$ws.Range("E11").Value = "  +3.95%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "'0.900"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "'15.04"
$ws.Range("E14").Value = "  +6.23%  "
$ws.Range("D15").Value = "2.341.08"
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("D16").Value = "'5.56"
$ws.Range("E16").Value = "  +6.10%  "
$ws.Range("D17").Value = "'20.52"
$ws.Range("E17").Value = "  +19.75%  "
$ws.Range("D18").Value = "2.043.75"
$ws.Range("E18").Value = "  +3.61%  "
$ws.Range("D19").Value = "37.558.18"
$ws.Range("E19").Value = "  +6.37%  "
$ws.Range("D20").Value = "'73.20"
$ws.Range("E20").Value = "  +4.76%  "
$ws.Range("D21").Value = "0.0₃0872"
$ws.Range("E21").Value = "  +4.52%  "
$ws.Range("D22").Value = "'5.33"
$ws.Range("E22").Value = "  +6.19%  "
$ws.Range("D23").Value = "'237.01"
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("D24").Value = "'2.68"
$ws.Range("E24").Value = "  +19.04%  "
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  +5.35%  "
$ws.Range("D27").Value = "'9.55"
$ws.Range("E27").Value = "  +5.59%  "
$ws.Range("D28").Value = "'164.53"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("D29").Value = "'19.87"
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("E31").Value = "  +8.96%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "'0.112"
$ws.Range("E32").Value = "  +26.14%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'1.21"
$ws.Range("E33").Value = "  +7.26%  "
$ws.Range("D34").Value = "'4.71"
$ws.Range("E34").Value = "  +10.86%  "
$ws.Range("D35").Value = "'0.0611"
$ws.Range("E35").Value = "  +4.92%  "
$ws.Range("D36").Value = "'2.44"
$ws.Range("E36").Value = "  +7.86%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("D39").Value = "'6.05"
$ws.Range("E39").Value = "  +25.04%  "
$ws.Range("E40").Value = "  +16.56%  "
$ws.Range("D41").Value = "'1.22"
$ws.Range("E41").Value = "  +3.80%  "
$ws.Range("E42").Value = "  +22.52%  "
$ws.Range("E43").Value = "  +2.80%  "
$ws.Range("D44").Value = "'0.0218"
$ws.Range("E44").Value = "  +5.17%  "
$ws.Range("D45").Value = "'1.13"
$ws.Range("E45").Value = "  +5.67%  "
$ws.Range("D46").Value = "'8.03"
$ws.Range("E46").Value = "  +8.47%  "
$ws.Range("D47").Value = "'16.81"
$ws.Range("E47").Value = "  +9.60%  "
$ws.Range("D48").Value = "'94.82"
$ws.Range("E48").Value = "  +5.02%  "
$ws.Range("D49").Value = "1.421.54"
$ws.Range("E49").Value = "  +4.18%  "
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("D51").Value = "'47.34"
$ws.Range("E51").Value = "  +4.39%  "
